# Update column F (dSF) values on specific rows to reflect repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F4").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("F12").Value = 11
$ws.Range("F14").Value = -3
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -4
$ws.Range("F19").Value = -4
$ws.Range("F23").Value = 5
